$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "admin"
$ws.Range("C3").Value = "admin"

$ws.Range("D3").Select()
